# Auto-generated edit script: update FFXIV leve profit market-price figures
# to match refreshed scheduled-runner market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 341.33334
$ws.Range("I33").Value = 273.66666
$ws.Range("J33").Value = 679.6667
$ws.Range("K33").Value = 273.66666
$ws.Range("L33").Value = 679.6667
$ws.Range("M33").Value = -44.66665999999998
$ws.Range("N33").Value = -1137.6667

# Row 43
$ws.Range("H43").Value = 1525
$ws.Range("I43").Value = 1150
$ws.Range("J43").Value = 1712.5
$ws.Range("K43").Value = 1150
$ws.Range("L43").Value = 1712.5
$ws.Range("M43").Value = -1081
$ws.Range("N43").Value = -1850.5

# Row 92
$ws.Range("H92").Value = 1373.5714
$ws.Range("I92").Value = 1373.5714
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1373.5714
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -125.5714
$ws.Range("N92").ClearContents()

# Row 98
$ws.Range("H98").Value = 2024.8055
$ws.Range("I98").Value = 310.52
$ws.Range("J98").Value = 5920.909
$ws.Range("K98").Value = 310.52
$ws.Range("L98").Value = 5920.909
$ws.Range("M98").Value = 1187.48
$ws.Range("N98").Value = -8916.909

# Row 103
$ws.Range("H103").Value = 1522.7778
$ws.Range("I103").Value = 1700
$ws.Range("J103").Value = 1301.25
$ws.Range("K103").Value = 5100
$ws.Range("L103").Value = 3903.75
$ws.Range("M103").Value = -4514
$ws.Range("N103").Value = -5075.75

# Row 122
$ws.Range("H122").Value = 2024.8055
$ws.Range("I122").Value = 310.52
$ws.Range("J122").Value = 5920.909
$ws.Range("K122").Value = 931.5599999999999
$ws.Range("L122").Value = 17762.727
$ws.Range("M122").Value = 1518.44
$ws.Range("N122").Value = -22662.727

# Row 129
$ws.Range("H129").Value = 991.3
$ws.Range("J129").Value = 997.2653
$ws.Range("L129").Value = 2991.7959
$ws.Range("N129").Value = -12991.7959

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 8633.615
$ws.Range("I61").Value = 9985.182000000001
$ws.Range("K61").Value = 9985.182000000001
$ws.Range("M61").Value = -9773.182000000001

# Row 74
$ws.Range("H74").Value = 6120.7
$ws.Range("I74").Value = 9222.666999999999
$ws.Range("J74").Value = 1467.75
$ws.Range("K74").Value = 9222.666999999999
$ws.Range("L74").Value = 1467.75
$ws.Range("M74").Value = -8348.666999999999
$ws.Range("N74").Value = -3215.75

# Row 77
$ws.Range("H77").Value = 6120.7
$ws.Range("I77").Value = 9222.666999999999
$ws.Range("J77").Value = 1467.75
$ws.Range("K77").Value = 46113.335
$ws.Range("L77").Value = 7338.75
$ws.Range("M77").Value = -41745.335
$ws.Range("N77").Value = -16074.75

# Row 129
$ws.Range("H129").Value = 39337.375
$ws.Range("J129").Value = 39337.375
$ws.Range("L129").Value = 39337.375
$ws.Range("N129").Value = -49337.375

# Row 136
$ws.Range("H136").Value = 8633.615
$ws.Range("I136").Value = 9985.182000000001
$ws.Range("K136").Value = 29955.546
$ws.Range("M136").Value = -27405.546

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2323.2307
$ws.Range("I20").Value = 2400.9
$ws.Range("J20").Value = 2064.3333
$ws.Range("K20").Value = 2400.9
$ws.Range("L20").Value = 2064.3333
$ws.Range("M20").Value = -2153.9
$ws.Range("N20").Value = -2558.3333

# Row 80
$ws.Range("H80").Value = 1846.125
$ws.Range("I80").Value = 1200
$ws.Range("J80").Value = 2139.818
$ws.Range("K80").Value = 1200
$ws.Range("L80").Value = 2139.818
$ws.Range("M80").Value = -202
$ws.Range("N80").Value = -4135.818

# Row 83
$ws.Range("H83").Value = 1846.125
$ws.Range("I83").Value = 1200
$ws.Range("J83").Value = 2139.818
$ws.Range("K83").Value = 6000
$ws.Range("L83").Value = 10699.09
$ws.Range("M83").Value = -1008
$ws.Range("N83").Value = -20683.09

# Row 96
$ws.Range("H96").Value = 11746.546
$ws.Range("I96").Value = 11746.546
$ws.Range("K96").Value = 11746.546
$ws.Range("M96").Value = -9000.546

# Row 105
$ws.Range("H105").Value = 2713.16
$ws.Range("I105").Value = 2761.625
$ws.Range("J105").Value = 1550
$ws.Range("K105").Value = 2761.625
$ws.Range("L105").Value = 1550
$ws.Range("M105").Value = -1014.625
$ws.Range("N105").Value = -5044

# Row 107
$ws.Range("H107").Value = 1650.1428
$ws.Range("I107").Value = 1162.2
$ws.Range("J107").Value = 2870
$ws.Range("K107").Value = 1162.2
$ws.Range("L107").Value = 2870
$ws.Range("M107").Value = 757.8
$ws.Range("N107").Value = -6710

# Row 140
$ws.Range("H140").Value = 89408
$ws.Range("J140").Value = 89408
$ws.Range("L140").Value = 89408
$ws.Range("N140").Value = -99768

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 969
$ws.Range("I22").Value = 1058.75
$ws.Range("J22").Value = 610
$ws.Range("K22").Value = 1058.75
$ws.Range("L22").Value = 610
$ws.Range("M22").Value = -708.75
$ws.Range("N22").Value = -1310

# Row 31
$ws.Range("H31").Value = 6190.0835
$ws.Range("I31").Value = 12047.111
$ws.Range("J31").Value = 2675.8667
$ws.Range("K31").Value = 12047.111
$ws.Range("L31").Value = 2675.8667
$ws.Range("M31").Value = -11752.111
$ws.Range("N31").Value = -3265.8667

# Row 33
$ws.Range("H33").Value = 1515.5
$ws.Range("I33").Value = 354
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 354
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = 25
$ws.Range("N33").Value = -5758

# Row 34
$ws.Range("H34").Value = 6190.0835
$ws.Range("I34").Value = 12047.111
$ws.Range("J34").Value = 2675.8667
$ws.Range("K34").Value = 12047.111
$ws.Range("L34").Value = 2675.8667
$ws.Range("M34").Value = -11845.111
$ws.Range("N34").Value = -3079.8667

# Row 58
$ws.Range("H58").Value = 967.1818
$ws.Range("I58").Value = 963.9
$ws.Range("K58").Value = 963.9
$ws.Range("M58").Value = -760.9

# Row 122
$ws.Range("H122").Value = 2565
$ws.Range("I122").Value = 2565
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7695
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5245
$ws.Range("N122").ClearContents()

# Row 134
$ws.Range("H134").Value = 10533.772
$ws.Range("I134").Value = 8530.357
$ws.Range("J134").Value = 14039.75
$ws.Range("K134").Value = 25591.071
$ws.Range("L134").Value = 42119.25
$ws.Range("M134").Value = -23056.071
$ws.Range("N134").Value = -47189.25

# Row 136
$ws.Range("H136").Value = 967.1818
$ws.Range("I136").Value = 963.9
$ws.Range("K136").Value = 2891.7
$ws.Range("M136").Value = -341.6999999999998

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 87.71429000000001
$ws.Range("I14").Value = 87.71429000000001
$ws.Range("K14").Value = 263.14287
$ws.Range("M14").Value = -90.14287000000002

# Row 68
$ws.Range("H68").Value = 26063
$ws.Range("I68").Value = 1124.5
$ws.Range("J68").Value = 51001.5
$ws.Range("K68").Value = 3373.5
$ws.Range("L68").Value = 153004.5
$ws.Range("M68").Value = -2562.5
$ws.Range("N68").Value = -154626.5

# Row 71
$ws.Range("H71").Value = 26063
$ws.Range("I71").Value = 1124.5
$ws.Range("J71").Value = 51001.5
$ws.Range("K71").Value = 10120.5
$ws.Range("L71").Value = 459013.5
$ws.Range("M71").Value = -6064.5
$ws.Range("N71").Value = -467125.5

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1153.36
$ws.Range("I102").Value = 1164.5454
$ws.Range("J102").Value = 1071.3334
$ws.Range("K102").Value = 1164.5454
$ws.Range("L102").Value = 1071.3334
$ws.Range("M102").Value = 457.4546
$ws.Range("N102").Value = -4315.3334

# Row 126
$ws.Range("H126").Value = 1018.4545
$ws.Range("I126").Value = 920.3
$ws.Range("K126").Value = 2760.9
$ws.Range("M126").Value = -290.8999999999996

# Row 132
$ws.Range("H132").Value = 19333.166
$ws.Range("I132").Value = 14625
$ws.Range("J132").Value = 28749.5
$ws.Range("K132").Value = 43875
$ws.Range("L132").Value = 86248.5
$ws.Range("M132").Value = -41345
$ws.Range("N132").Value = -91308.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 618.4545000000001
$ws.Range("I22").Value = 498.2
$ws.Range("J22").Value = 718.6667
$ws.Range("K22").Value = 498.2
$ws.Range("L22").Value = 718.6667
$ws.Range("M22").Value = -203.2
$ws.Range("N22").Value = -1308.6667

# Row 27
$ws.Range("H27").Value = 618.4545000000001
$ws.Range("I27").Value = 498.2
$ws.Range("J27").Value = 718.6667
$ws.Range("K27").Value = 498.2
$ws.Range("L27").Value = 718.6667
$ws.Range("M27").Value = -391.2
$ws.Range("N27").Value = -932.6667
